# Updates the cryptos price/volume snapshot (GitHub Actions data refresh).
# Numeric-looking "Price" values are written with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.387.61"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "2.949.70"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'495.77"
$ws.Range("E5").Value = "  -5.92%  "
$ws.Range("D6").Value = "'134.07"
$ws.Range("E6").Value = "  -6.38%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -5.35%  "
$ws.Range("E9").Value = "  -6.32%  "
$ws.Range("E10").Value = "  -6.50%  "
$ws.Range("D11").Value = "'0.352"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("D12").Value = "3.460.24"
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").Value = "'0.125"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").Value = "'25.87"
$ws.Range("E14").Value = "  -5.29%  "
$ws.Range("E15").Value = "  -8.48%  "
$ws.Range("D16").Value = "56.548.92"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("D18").Value = "2.947.73"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("E19").Value = "  -5.65%  "
$ws.Range("D20").Value = "'7.74"
$ws.Range("E20").Value = "  -5.70%  "
$ws.Range("D21").Value = "'315.85"
$ws.Range("E21").Value = "  -7.47%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").Value = "'62.35"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -5.65%  "
$ws.Range("E28").Value = "  -11.11%  "
$ws.Range("D29").Value = "'6.46"
$ws.Range("E29").Value = "  -7.38%  "
$ws.Range("E30").Value = "  -6.09%  "
$ws.Range("E31").Value = "  -6.36%  "
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("E33").Value = "  -8.90%  "
$ws.Range("D34").Value = "'152.51"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  -6.65%  "
$ws.Range("E36").Value = "  -5.41%  "
$ws.Range("D37").Value = "'1.20"
$ws.Range("E37").Value = "  -9.58%  "
$ws.Range("D38").Value = "'23.64"
$ws.Range("E38").Value = "  -9.97%  "
$ws.Range("D39").Value = "'0.0653"
$ws.Range("E39").Value = "  -7.31%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "2.981.16"
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'37.30"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D45").Value = "2.143.72"
$ws.Range("E45").Value = "  -8.45%  "
$ws.Range("E46").Value = "  -8.73%  "
$ws.Range("D47").Value = "'5.84"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").Value = "'0.915"
$ws.Range("E48").Value = "  -11.48%  "
$ws.Range("D49").Value = "'0.0230"
$ws.Range("E49").Value = "  -5.97%  "
$ws.Range("E50").Value = "  -6.15%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0856"
$ws.Range("E51").Value = "  -5.12%  "
